# TaskList.xlsx update:
#  - Task 1 ("1. Reduce count of Gomel-sat news to 45 (3 pages)") status
#    moves from "In progress" to "Done", highlighted in green.
#  - Task 5 ("5. Search algorithm optimization") status moves from "Open"
#    to "In Progress" and gets an assignee (Arthur).
#  - A brand-new task 13 ("13. Fix bug of incorrect div count") is added
#    in the next empty row, with Priority "Hight" and Status "Open".
#  - Selection cursor ends up on E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 1: status "In progress" -> "Done", shown in green.
$ws.Range("C2").Value = "Done"
$ws.Range("C2").Font.Color = 5287936   # RGB(0, 176, 80) == FF00B050

# New task 13, placed in the next free row of the list.
$ws.Range("A14").Value = "13. Fix bug of incorrect div count"
$ws.Range("B14").Value = "Hight"
$ws.Range("C14").Value = "Open"

# Task 5: status "Open" -> "In Progress", assignee "Arthur".
$ws.Range("C6").Value = "In Progress"
$ws.Range("D6").Value = "Arthur"

# Update the active selection.
$ws.Range("E9").Select()
